$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.000.23'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').Value = '3.151.44'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('E4').Value = '  -0.08%  '
$c = $ws.Range('D5')
$c.Value = '''573.84'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +2.59%  '
$c = $ws.Range('D6')
$c.Value = '''149.75'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +5.07%  '
$c = $ws.Range('D7')
$c.Value = '''0.999'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.149.98'
$ws.Range('E8').Value = '  +3.00%  '
$c = $ws.Range('D9')
$c.Value = '''0.526'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +1.94%  '
$c = $ws.Range('D10')
$c.Value = '''0.160'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +3.99%  '
$c = $ws.Range('D11')
$c.Value = '''6.15'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +0.15%  '
$c = $ws.Range('D12')
$c.Value = '''0.498'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +3.63%  '
$c = $ws.Range('D13')
$c.Value = '''0.0000263'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +13.24%  '
$c = $ws.Range('D14')
$c.Value = '''37.06'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +4.90%  '
$ws.Range('D15').Value = '3.665.45'
$ws.Range('E15').Value = '  +2.90%  '
$ws.Range('D16').Value = '65.077.16'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('D17').Value = '3.149.92'
$ws.Range('E17').Value = '  +2.75%  '
$c = $ws.Range('D18')
$c.Value = '''7.10'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +4.69%  '
$ws.Range('E19').Value = '  +1.17%  '
$c = $ws.Range('D20')
$c.Value = '''507.21'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +4.26%  '
$c = $ws.Range('D21')
$c.Value = '''14.77'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +3.17%  '
$c = $ws.Range('D22')
$c.Value = '''0.717'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +3.89%  '
$c = $ws.Range('D23')
$c.Value = '''15.27'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +4.05%  '
$c = $ws.Range('D24')
$c.Value = '''7.72'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +2.66%  '
$c = $ws.Range('D25')
$c.Value = '''84.19'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('E27').Value = '  +3.49%  '
$ws.Range('E28').Value = '  +8.06%  '
$c = $ws.Range('D29')
$c.Value = '''2.17'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +5.71%  '
$c = $ws.Range('D30')
$c.Value = '''2.80'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +8.61%  '
$c = $ws.Range('D31')
$c.Value = '''27.60'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +4.18%  '
$ws.Range('E32').Value = '  -0.13%  '
$c = $ws.Range('D33')
$c.Value = '''1.19'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +2.24%  '
$c = $ws.Range('D34')
$c.Value = '''6.20'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +8.90%  '
$c = $ws.Range('D35')
$c.Value = '''6.52'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +4.49%  '
$c = $ws.Range('D36')
$c.Value = '''54.97'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.28%  '
$c = $ws.Range('D37')
$c.Value = '''0.0897'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +10.23%  '
$c = $ws.Range('D38')
$c.Value = '''464.20'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +4.41%  '
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('E40').Value = '  +7.77%  '
$c = $ws.Range('D41')
$c.Value = '''8.66'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +3.91%  '
$ws.Range('D42').Value = '3.051.79'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('E43').Value = '  +0.02%  '
$c = $ws.Range('D44')
$c.Value = '''2.44'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +8.13%  '
$ws.Range('E45').Value = '  +2.18%  '
$c = $ws.Range('D46')
$c.Value = '''28.55'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('D47').Value = '0.0₃0583'
$ws.Range('E47').Value = '  +12.37%  '
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('E50').Value = '  +5.17%  '
$c = $ws.Range('D51')
$c.Value = '''119.43'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +1.36%  '
